$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-10 from
# 2023-10-05 (serial 45204) to 2023-10-08 (serial 45207).
# The cells are date-formatted, so reading .Value() back yields a
# [datetime], not a plain number - compare against a converted datetime.
$oldDate = [datetime]::FromOADate(45204)

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = 45207
    }
}
